$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new log entry row (row 4) to the BBM upload log table
$ws.Range("A4").Value = "KIRAN KUMAR"
$ws.Range("B4").Value = "OS"
$ws.Range("C4").Value = "Ftth OS_25.11.2025.xlsx"
$ws.Range("D4").Value = "2025-12-02 11:47"
$ws.Range("E4").Value = "2025-12"
